$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the current newest-date header (currently in column B) before
# the columns shift, so it can be carried into its new position.
$previousHeader = $ws.Range("B1").Value2

# Insert two new columns before the existing column C.
# This shifts the old column C (dated "Jun_10" header) to column E.
$ws.Columns("C:D").Insert()

# The old B1 header value now belongs in D1 (it shifted two columns right
# conceptually, but Insert() only carried column C's data to E, so restore
# the old B1 value into D1 explicitly).
$ws.Range("D1").Value = $previousHeader

# New header row: C1 becomes the next date header ("Jun_15"), and B1
# becomes the newest date header ("Jun_17"). E1 already carries the old
# C1 header value ("Jun_10") because of the column insert.
$ws.Range("C1").Value = "Jun_15"
$ws.Range("B1").Value = "Jun_17"

# Fill the two new data columns (C and D) with the same default rating
# value ("UN") that already fills column B for every data row.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# Restore explicit custom width (8.0) on the three right-most columns:
# the two newly inserted columns (C, D) and the shifted former-C column
# now sitting at E.
$ws.Columns("C").ColumnWidth = 7.14
$ws.Columns("D").ColumnWidth = 7.14
$ws.Columns("E").ColumnWidth = 7.14
